$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packing Slip")

# Fill G12:G15 with the same relative formula as G11 (F*E), creating a shared formula group
$ws.Range("G12:G15").Formula = "=F12*E12"

# Label + subtotal row
$ws.Range("F16").Value = "Subtotal"
$ws.Range("G16").Formula = "=SUM(G11:G15)"

# Move the active selection to G16, matching the author's final cursor position
$ws.Range("G16").Select()
